$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C1 now points at the only remaining "project1" shared string ---
$ws.Range("C1").Value = "project1"

# --- Rows 2 and 3 lose their data (B/C, and A's text) but keep their
#     existing "Hyperlink" cell style (ClearContents leaves formatting
#     untouched); row 4 is removed entirely ---
$ws.Range("A2:C3").ClearContents()
$ws.Rows("4:4").Delete()

# --- Remove the hyperlinks on (the former) A2, A3 and A4, keeping the one
#     on A1. Hyperlinks.Delete() on a scoped Range wipes the whole sheet's
#     collection in this host, so collect the ones to remove first and
#     delete them in reverse index order to avoid the live-collection
#     reindexing problem. ---
$toRemove = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -ne "`$A`$1") {
        $toRemove += $h
    }
}
for ($i = $toRemove.Count - 1; $i -ge 0; $i--) {
    $toRemove[$i].Delete()
}

# --- Selection moves to A2 ---
$ws.Range("A2").Select()
